$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8 to make room for the new article (shifts old rows 8-17 down to 9-18)
$ws.Rows.Item(8).Insert()

$ws.Range("A7").Value = "...アンソロピックなど新興企業がAIサービスを発表すると、代替されるとの懸念からフェーズ4のソフトウエアや保険、資産運用サービスの一角が急落した。
                          　..."
$ws.Range("B7").Value = "https://www.nikkei.com/article/DGXZQOUB09A630Z00C26A2000000/"
$ws.Range("A8").Value = "...ヌビーンは米国教職員退職年金保険組合（TIAA）の運用部門で、同社が新たに設立した子会社を通じてシュローダーの全株式を取得する。買収提案は1株あたり6.12ポンド（配当含む）と11日の終値に対し34％のプレミアム（上乗せ幅）を乗せた水準だ。
                          　..."
$ws.Range("B8").Value = "https://www.nikkei.com/article/DGXZQOUB127BB0S6A210C2000000/"
$ws.Range("A9").Value = "...設備投資は英政府が25年4月に企業の国民保険料負担を引き上げたことで低迷している。
                          　25年通年のGDPの伸び率は1.3%だった。24年の1.1%をわずかに上回った。..."
$ws.Range("B9").Value = "https://www.nikkei.com/article/DGXZQOGR1262F0S6A210C2000000/"
$ws.Range("A10").Value = "マイナ保険証の利用率、12月に6割超え　紙の保険証の期限切れ"
$ws.Range("B10").Value = "https://www.nikkei.com/article/DGXZQOUA123B50S6A210C2000000/"
$ws.Range("A11").Value = "...
                          　厚生労働省は12日、マイナンバーカードと健康保険証を一体にした「マイナ保険証」の2025年12月の利用率が63.24%だったと明らかにした。..."
$ws.Range("B11").Value = "https://www.nikkei.com/article/DGXZQOUA123B50S6A210C2000000/"
$ws.Range("A12").Value = "...三井住友海上火災保険も25年11月に欧米で資産運用事業を手掛ける米ベアリングスに14億4000万ドル（約2200億円）出資すると発表した。みずほ銀行は24年に「プライベートデット（非上場企業向けの融資）」に強みを持つ米資産運用会社のゴラブ・キャピタルに出資している。..."
$ws.Range("B12").Value = "https://www.nikkei.com/article/DGXZQOUB101P00Q6A210C2000000/"
$ws.Range("A13").Value = "...安野貴博党首は11日に公開したユーチューブ配信で、消費税率の維持や社会保険料の引き下げを訴えたのを勝因の一つだと分析した。「ほぼ全ての政党が消費税の減税を掲げるなかで、社会保険料の引き下げに支持をいただいた部分も一定あった」と語った。..."
$ws.Range("B13").Value = "https://www.nikkei.com/article/DGXZQOUA1119K0R10C26A2000000/"
$ws.Range("A14").Value = "...代理店への出向者や傘下の生命保険会社の担当部門を対象に、アンケート・ヒアリング調査のほか、電子メールなどのデジタルフォレンジック（電子鑑識）や会議資料の確認を実施した。
                          　..."
$ws.Range("B14").Value = "https://www.nikkei.com/article/DGXZQOUB1212T0S6A210C2000000/"
$ws.Range("A15").Value = "米国の失業保険申請、週22.7万件　市場予想上回る"
$ws.Range("B15").Value = "https://www.nikkei.com/article/DGXZQOGN1149Z0R10C26A2000000/"
$ws.Range("A16").Value = "...
                          　【ワシントン=野一色遥花】米労働省が12日発表した失業保険統計（季節調整済み）によると、企業の解雇動向を映す失業保険の新規申請件数は2月1〜7日に22万7000件だった。..."
$ws.Range("B16").Value = "https://www.nikkei.com/article/DGXZQOGN1149Z0R10C26A2000000/"
$ws.Range("A17").Value = "...12日発表の週間の新規失業保険申請件数は22万7000件と市場予想をやや上回ったが、前の週からは減った。
                          　市場では「労働市場の環境は安定しており、2026年に入って徐々に改善する兆しをみせている」（オックスフォード・エコノミクスのナンシー・ホウテン氏）との受け止めがある。..."
$ws.Range("B17").Value = "https://www.nikkei.com/article/DGXZQOFL12BRU0S6A210C2000000/"
$ws.Range("A18").Value = "...12日発表の週間の新規失業保険申請件数は22万7000件と市場予想をやや上回ったが、前の週からは減った。
                          　市場では「労働市場の環境は安定しており、2026年に入って徐々に改善する兆しをみせている」（オックスフォード・エコノミクスのナンシー・ホウテン氏）との受け止めがある。..."
$ws.Range("B18").Value = "https://www.nikkei.com/article/DGXZQOFL12BFD0S6A210C2000000/"

# The newly inserted row 8 needs date/category filled in (other rows already carry these over from the shift)
# Leading apostrophe forces the date-shaped string to stay text (matches the other rows' storage as text, not a date value)
$ws.Range("C8").Value = "'2026-02-12"
$ws.Range("D8").Value = "最新動態"

